# Fixed lost data through imports
# Adds a new "ON-STREAM" technology entry (column D) on the
# "Object Type 2" sheet, which had been dropped by a prior import.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Object Type 2")
$ws.Activate()

# Row 12 (H2 System / Bus) previously only had object_index/parameter_index
# up to column C ("Technology"); restore the missing "Year" parameter_index
# and shift the "Technology" value into the new column D.
$ws.Range("C12").Value = "Year"
$ws.Range("D12").Value = "Technology"

# Rows 26-28 each gain a 4th parameter_index value that was lost on import.
$ws.Range("D26").Value = "ON-STREAM"
$ws.Range("D27").Value = "Bus"
$ws.Range("D28").Value = "Time"

$ws2 = $wb.Worksheets.Item("Relationship Type 1")
$ws2.Activate()
$ws2.Range("D56:D57").Select()

$ws.Activate()
$ws.Range("D28").Select()
